$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.138.50"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.797.11"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.98"
$ws.Range("D5").Style = $ws.Range("B2").Style
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5198"
$ws.Range("D7").Style = $ws.Range("B2").Style
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3818"
$ws.Range("D8").Style = $ws.Range("B2").Style
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07895"
$ws.Range("D9").Style = $ws.Range("B2").Style
$ws.Range("E9").Value = "  -4.31%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.100"
$ws.Range("D11").Style = $ws.Range("B2").Style
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.284"
$ws.Range("D12").Style = $ws.Range("B2").Style
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.56"
$ws.Range("D14").Style = $ws.Range("B2").Style
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "1.790.04"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.241"
$ws.Range("D16").Style = $ws.Range("B2").Style
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.89"
$ws.Range("D17").Style = $ws.Range("B2").Style
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001084"
$ws.Range("D18").Style = $ws.Range("B2").Style
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06556"
$ws.Range("D19").Style = $ws.Range("B2").Style
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.27"
$ws.Range("D21").Style = $ws.Range("B2").Style
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.950"
$ws.Range("D22").Style = $ws.Range("B2").Style
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").Value = "28.166.74"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("D24").Style = $ws.Range("B2").Style
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.270"
$ws.Range("D25").Style = $ws.Range("B2").Style
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.38"
$ws.Range("D26").Style = $ws.Range("B2").Style
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.44"
$ws.Range("D27").Style = $ws.Range("B2").Style
$ws.Range("E27").Value = "  -4.37%  "
$ws.Range("D28").Value = "1.999.98"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.331"
$ws.Range("D29").Style = $ws.Range("B2").Style
$ws.Range("E29").Value = "  -3.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.84"
$ws.Range("D30").Style = $ws.Range("B2").Style
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1063"
$ws.Range("D31").Style = $ws.Range("B2").Style
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.050"
$ws.Range("D32").Style = $ws.Range("B2").Style
$ws.Range("E32").Value = "  -5.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.677"
$ws.Range("D33").Style = $ws.Range("B2").Style
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.565"
$ws.Range("D34").Style = $ws.Range("B2").Style
$ws.Range("E34").Value = "  -3.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07327"
$ws.Range("D35").Style = $ws.Range("B2").Style
$ws.Range("E35").Value = "  +3.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.20"
$ws.Range("D36").Style = $ws.Range("B2").Style
$ws.Range("E36").Value = "  +8.11%  "
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2140"
$ws.Range("D38").Style = $ws.Range("B2").Style
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.714"
$ws.Range("D39").Style = $ws.Range("B2").Style
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.068"
$ws.Range("D40").Style = $ws.Range("B2").Style
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6143"
$ws.Range("D41").Style = $ws.Range("B2").Style
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.159"
$ws.Range("D42").Style = $ws.Range("B2").Style
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.373"
$ws.Range("D43").Style = $ws.Range("B2").Style
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.28"
$ws.Range("D44").Style = $ws.Range("B2").Style
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.778"
$ws.Range("D45").Style = $ws.Range("B2").Style
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5954"
$ws.Range("D46").Style = $ws.Range("B2").Style
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.79"
$ws.Range("D47").Style = $ws.Range("B2").Style
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.232"
$ws.Range("D48").Style = $ws.Range("B2").Style
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.915"
$ws.Range("D49").Style = $ws.Range("B2").Style
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06766"
$ws.Range("D50").Style = $ws.Range("B2").Style
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.15"
$ws.Range("D51").Style = $ws.Range("B2").Style
$ws.Range("E51").Value = "  -1.55%  "
